$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("segments"), shifting the old
# B:K block (PercActivations .. totalStd) one column to the right (C:L).
$ws.Columns("B").Insert()

# The inserted column B picked up column A's header formatting; copy the
# (now shifted) header style from C1 so B1 matches the other header cells.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# The data cells in the new column B also inherited column A's bold/
# bordered style from the insert - the segment-name values should be
# plain (unstyled), like the rest of the data columns.
$ws.Range("B2:B20").Style = "Normal"

# Rows 3-20 never had totalMean/totalStd (only row 2 does); the insert
# turned those originally-blank cells (old J3:K20) into real empty
# strings at K3:L20 - clear them back out so they stay genuinely blank.
$ws.Range("K3:L20").ClearContents()

# Column A becomes a 0-based numeric index per segment row, and column B
# receives the segment names that used to live in column A.
$segments = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

for ($i = 0; $i -lt $segments.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $segments[$i]
}
